# Update Pr(>F) p-values in column F (rounded to 3 digits after re-run
# without the log(x+1) y-axis scale on the biomass data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 0.0105
    4  = 0.4676
    7  = 0.1666
    8  = 0.0511
    9  = 0.9552
    12 = 0.346
    13 = 0.5948
    14 = 0.4058
    17 = 0.0729
    18 = 0.9962
    19 = 0.5725
    22 = 0.1733
    23 = 0.4874
    24 = 0.2239
    27 = 0.8888
    28 = 0.0143
    29 = 0.8919
    32 = 0.1679
    33 = 0.7092
    34 = 0.9668
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
